$d = $word.ActiveDocument

$spaces = " " * 66
$oldText = $spaces + "S.Gayathri"

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $oldText, 0)
